# Append " - Bitcoin" (as its own run, matching the title's bold/56pt
# formatting) right after the existing " Research" run in the title
# paragraph, e.g. "Network Protocol Research" -> "Network Protocol
# Research - Bitcoin".
#
# A plain InsertAfter()/TypeText() on the title paragraph creates a new
# run but with no run-properties (rPr) at all, and a Find/Replace or
# Range.Text= edit that touches the existing " Research" run causes it
# to be re-merged with its neighbours, losing the "separate run" shape
# that the target document has. So instead: build the new run's text in
# a scratch paragraph (where edits are free to coalesce however they
# like), seeded via Copy/Paste so it inherits the exact formatting
# (rFonts/b/bCs/sz/szCs) of the existing title run, then Copy/Paste the
# finished run into place and remove the scratch paragraph.

$d = $word.ActiveDocument

$title = $d.Paragraphs(2)
$titleRange = $title.Range
$insertPoint = $titleRange.End - 1   # just after "Research", before the pilcrow

# 1) Scratch paragraph at the very end of the document.
$docEnd = $d.Content.End
$d.Range($docEnd, $docEnd).InsertParagraphAfter()
$scratch = $d.Paragraphs($d.Paragraphs.Count)
$scratchStart = $scratch.Range.Start

# 2) Seed the scratch paragraph with a correctly-formatted copy of
#    " Research" (same run as what precedes our insertion point).
$sampleRange = $d.Range($insertPoint - 9, $insertPoint)
$sampleRange.Copy()
$d.Range($scratchStart, $scratchStart).Paste()

# 3) Turn the scratch paragraph's text from " Research" into
#    " - Bitcoin", keeping it a single run (formatting carries over).
$d.Range($scratchStart + 1, $scratchStart + 9).Delete()
$d.Range($scratchStart, $scratchStart + 1).InsertAfter("- Bitcoin")

# 4) Copy the finished " - Bitcoin" run and paste it right after
#    " Research" in the title paragraph.
$finalSrc = $d.Range($scratchStart, $scratchStart + 10)
$finalSrc.Copy()
$d.Range($insertPoint, $insertPoint).Paste()

# 5) Remove the scratch paragraph (now shifted to the new end of doc).
$scratch2 = $d.Paragraphs($d.Paragraphs.Count)
$scratch2.Range.Delete()

Write-Output $title.Range.Text
